$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 4225
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 4225
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 12675
$ws.Range("M17").ClearContents()
$ws.Range("N17").Value = -13011

# Row 32
$ws.Range("H32").Value = 1359.6
$ws.Range("I32").Value = 1899.5
$ws.Range("K32").Value = 1899.5
$ws.Range("M32").Value = -1573.5

# Row 33
$ws.Range("H33").Value = 282.4375
$ws.Range("I33").Value = 193.92308
$ws.Range("K33").Value = 193.92308
$ws.Range("M33").Value = 35.07692

# Row 86
$ws.Range("H86").Value = 2182.2942
$ws.Range("I86").Value = 2000.7
$ws.Range("J86").Value = 2441.7144
$ws.Range("K86").Value = 2000.7
$ws.Range("L86").Value = 2441.7144
$ws.Range("M86").Value = -877.7
$ws.Range("N86").Value = -4687.7144

# Row 89
$ws.Range("H89").Value = 2182.2942
$ws.Range("I89").Value = 2000.7
$ws.Range("J89").Value = 2441.7144
$ws.Range("K89").Value = 10003.5
$ws.Range("L89").Value = 12208.572
$ws.Range("M89").Value = -4387.5
$ws.Range("N89").Value = -23440.572

# Row 112
$ws.Range("H112").Value = 2271.8928
$ws.Range("J112").Value = 2535.389
$ws.Range("L112").Value = 7606.167
$ws.Range("N112").Value = -9822.167000000001

# Row 113
$ws.Range("H113").Value = 4000
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 4000
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 4000
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -10508

# Row 115
$ws.Range("H115").Value = 612.125

# Row 129
$ws.Range("H129").Value = 388427.06
$ws.Range("I129").Value = 410981.6
$ws.Range("K129").Value = 1232944.8
$ws.Range("M129").Value = -1227944.8

# Row 138
$ws.Range("H138").Value = 2115.9275
$ws.Range("I138").Value = 780.06665
$ws.Range("J138").Value = 2487
$ws.Range("K138").Value = 2340.19995
$ws.Range("L138").Value = 7461
$ws.Range("M138").Value = 2799.80005
$ws.Range("N138").Value = -17741

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 127210.51
$ws.Range("I32").Value = 139617.17
$ws.Range("K32").Value = 139617.17
$ws.Range("M32").Value = -139330.17

# Row 45
$ws.Range("H45").Value = 2117.9092
$ws.Range("I45").Value = 2216.1667
$ws.Range("K45").Value = 2216.1667
$ws.Range("M45").Value = -1839.1667

# Row 74
$ws.Range("H74").Value = 2532108.8
$ws.Range("I74").Value = 3474179
$ws.Range("K74").Value = 3474179
$ws.Range("M74").Value = -3473305

# Row 77
$ws.Range("H77").Value = 2532108.8
$ws.Range("I77").Value = 3474179
$ws.Range("K77").Value = 17370895
$ws.Range("M77").Value = -17366527

# Row 132
$ws.Range("H132").Value = 642377.1
$ws.Range("I132").Value = 695555.75
$ws.Range("K132").Value = 2086667.25
$ws.Range("M132").Value = -2084137.25

$ws = $wb.Worksheets.Item("CRP")
# Row 132
$ws.Range("H132").Value = 2478.65
$ws.Range("I132").Value = 2121.1765
$ws.Range("J132").Value = 4504.3335
$ws.Range("K132").Value = 6363.529500000001
$ws.Range("L132").Value = 13513.0005
$ws.Range("M132").Value = -3833.529500000001
$ws.Range("N132").Value = -18573.0005

$ws = $wb.Worksheets.Item("CUL")
# Row 87
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()

# Row 88
$ws.Range("H88").Value = 48571.145
$ws.Range("I88").Value = 29999
$ws.Range("J88").Value = 51666.5
$ws.Range("K88").Value = 89997
$ws.Range("L88").Value = 154999.5
$ws.Range("M88").Value = -89569
$ws.Range("N88").Value = -155855.5

# Row 90
$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()

# Row 91
$ws.Range("H91").Value = 48571.145
$ws.Range("I91").Value = 29999
$ws.Range("J91").Value = 51666.5
$ws.Range("K91").Value = 89997
$ws.Range("L91").Value = 154999.5
$ws.Range("M91").Value = -88515
$ws.Range("N91").Value = -157963.5

# Row 107
$ws.Range("H107").Value = 1940
$ws.Range("I107").Value = 881.5
$ws.Range("J107").Value = 2998.5
$ws.Range("K107").Value = 2644.5
$ws.Range("L107").Value = 8995.5
$ws.Range("M107").Value = -724.5
$ws.Range("N107").Value = -12835.5

# Row 109
$ws.Range("H109").Value = 1650.2667
$ws.Range("I109").Value = 903.4
$ws.Range("J109").Value = 3144
$ws.Range("K109").Value = 2710.2
$ws.Range("L109").Value = 9432
$ws.Range("M109").Value = -1670.2
$ws.Range("N109").Value = -11512

$ws = $wb.Worksheets.Item("GSM")
# Row 113
$ws.Range("H113").Value = 2976.0386
$ws.Range("I113").Value = 2170.7646
$ws.Range("K113").Value = 2170.7646
$ws.Range("M113").Value = -0.7645999999999731

$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 4923.6924
$ws.Range("I40").Value = 4592.1816
$ws.Range("J40").Value = 6747
$ws.Range("K40").Value = 4592.1816
$ws.Range("L40").Value = 6747
$ws.Range("M40").Value = -4456.1816
$ws.Range("N40").Value = -7019

# Row 131
$ws.Range("H131").Value = 50000
$ws.Range("J131").Value = 50000
$ws.Range("L131").Value = 50000
$ws.Range("N131").Value = -60080

# Row 132
$ws.Range("H132").Value = 8349374.5
$ws.Range("I132").Value = 14608681
$ws.Range("J132").Value = 3632.5
$ws.Range("K132").Value = 43826043
$ws.Range("L132").Value = 10897.5
$ws.Range("M132").Value = -43823513
$ws.Range("N132").Value = -15957.5

# Row 136
$ws.Range("H136").Value = 10872410
$ws.Range("I136").Value = 12502338
$ws.Range("K136").Value = 37507014
$ws.Range("M136").Value = -37504464

$ws = $wb.Worksheets.Item("WVR")
# Row 69
$ws.Range("H69").Value = 72499
$ws.Range("J69").Value = 72499
$ws.Range("L69").Value = 72499
$ws.Range("N69").Value = -73997

# Row 72
$ws.Range("H72").Value = 72499
$ws.Range("J72").Value = 72499
$ws.Range("L72").Value = 217497
$ws.Range("N72").Value = -224985

# Row 107
$ws.Range("H107").Value = 1917.375
$ws.Range("I107").Value = 1092.2106
$ws.Range("J107").Value = 5053
$ws.Range("K107").Value = 3276.6318
$ws.Range("L107").Value = 15159
$ws.Range("M107").Value = -1356.6318
$ws.Range("N107").Value = -18999

# Row 113
$ws.Range("H113").Value = 572.73334
$ws.Range("I113").Value = 651
$ws.Range("J113").Value = 437.54544
$ws.Range("K113").Value = 1953
$ws.Range("L113").Value = 1312.63632
$ws.Range("M113").Value = 217
$ws.Range("N113").Value = -5652.63632

# Row 119
$ws.Range("H119").Value = 69997
$ws.Range("J119").Value = 69997
$ws.Range("L119").Value = 69997
$ws.Range("N119").Value = -79673

# Row 122
$ws.Range("H122").Value = 37452.16
$ws.Range("I122").Value = 1073.1154
$ws.Range("J122").Value = 226623.2
$ws.Range("K122").Value = 3219.3462
$ws.Range("L122").Value = 679869.6000000001
$ws.Range("M122").Value = -769.3462
$ws.Range("N122").Value = -684769.6000000001

# Row 132
$ws.Range("H132").Value = 9806261
$ws.Range("I132").Value = 12822550
$ws.Range("J132").Value = 3320.75
$ws.Range("K132").Value = 38467650
$ws.Range("L132").Value = 9962.25
$ws.Range("M132").Value = -38465120
$ws.Range("N132").Value = -15022.25

# Row 136
$ws.Range("H136").Value = 17314502
$ws.Range("I136").Value = 8597872
$ws.Range("K136").Value = 25793616
$ws.Range("M136").Value = -25791066
